# Generate Report for Handoff
#
# Updates the "localization-status" workbook with a new handoff run:
#  - the source markdown file got a new guid/name and a later handoff timestamp
#  - two new source files (.png) were picked up and handed off as well
#
# Sheet 1 = "Overview" (cols A:D)
# Sheet 2 = "zh-cn"    (cols A:L)
# Sheet 3 = "de-de"    (cols A:L)

$wb = $excel.ActiveWorkbook

$oldMd  = "8a4d39d5-7a86-4cbc-b43b-540740e36335.md"
$newMd  = "21f3b798-25d7-4e08-ad42-69f73d4ca63d.md"

$newXlfZh = "21f3b798-25d7-4e08-ad42-69f73d4ca63d.5636a32b385bfce3db27335ccde3d6281d4c3c3c.zh-cn.xlf"
$newXlfDe = "21f3b798-25d7-4e08-ad42-69f73d4ca63d.5636a32b385bfce3db27335ccde3d6281d4c3c3c.de-de.xlf"

$png1 = "27204dba-f0f4-4b63-8d76-870c65938ac6.png"
$png2 = "91eed04d-89c8-46bb-9416-8bf3d53fa8a2.png"

$pngTargetZh1 = "3fe3491d92f1a28aee782048bfab3f2e8fc1dca6.png"
$pngTargetZh2 = "73534a6839b4b303167415f968c674537981d2dc.png"
$pngTargetDe1 = "3fe3491d92f1a28aee782048bfab3f2e8fc1dca6.png"
$pngTargetDe2 = "73534a6839b4b303167415f968c674537981d2dc.png"

$overviewDate = "2016-03-23 19:15:21"
$zhHandoffDate = "2016-03-23 19:15:16"
$deHandoffDate = "2016-03-23 19:15:21"
$epoch = "0001-01-01 00:00:00"
$status = "Ready for handoff"
$dateFmt = "yyyy-mm-dd HH:mm:ss"

$mdBase  = "https://github.com/OpenLocalizationTest/oltest/blob/f42d5a4c3e66cdfdc0c84a26f9cbe5e5494c93fb/e2e"
$zhBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b88569b97e266ace5684d4c08ab8fb4411cbe8c7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2140457273b0a8ff43fd39694b6aee3258571e27/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

# =====================================================================
# Sheet 1: "Overview"
# =====================================================================
$ws1 = $wb.Worksheets.Item(1)
$ws1.Hyperlinks.Delete()

$ws1.Cells.Item(2,1).Value = $newMd
$ws1.Cells.Item(2,2).Value = $status
$ws1.Cells.Item(2,3).Value = $status
$ws1.Cells.Item(2,4).Value = $overviewDate
$ws1.Cells.Item(2,4).NumberFormat = $dateFmt

$ws1.Cells.Item(3,1).Value = $png1
$ws1.Cells.Item(3,2).Value = $status
$ws1.Cells.Item(3,3).Value = $status
$ws1.Cells.Item(3,4).Value = $overviewDate
$ws1.Cells.Item(3,4).NumberFormat = $dateFmt

$ws1.Cells.Item(4,1).Value = $png2
$ws1.Cells.Item(4,2).Value = $status
$ws1.Cells.Item(4,3).Value = $status
$ws1.Cells.Item(4,4).Value = $overviewDate
$ws1.Cells.Item(4,4).NumberFormat = $dateFmt

$ws1.Hyperlinks.Add($ws1.Range("A2"), "$mdBase/$newMd", "", "", $newMd) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "$mdBase/$png1", "", "", $png1) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "$mdBase/$png2", "", "", $png2) | Out-Null

# =====================================================================
# Sheet 2: "zh-cn"
# =====================================================================
$ws2 = $wb.Worksheets.Item(2)
$ws2.Hyperlinks.Delete()

$ws2.Cells.Item(2,1).Value = $newMd
$ws2.Cells.Item(2,4).Value = $newXlfZh
$ws2.Cells.Item(2,5).Value = $zhHandoffDate
$ws2.Cells.Item(2,5).NumberFormat = $dateFmt

$ws2.Cells.Item(3,1).Value = $png1
$ws2.Cells.Item(3,2).Value = ".png"
$ws2.Cells.Item(3,3).Value = $status
$ws2.Cells.Item(3,4).Value = $pngTargetZh1
$ws2.Cells.Item(3,5).Value = $zhHandoffDate
$ws2.Cells.Item(3,5).NumberFormat = $dateFmt
$ws2.Cells.Item(3,8).Value = $epoch
$ws2.Cells.Item(3,8).NumberFormat = $dateFmt
$ws2.Cells.Item(3,10).Value = "IsDependency"
$ws2.Cells.Item(3,11).Value = "e2e\$newMd"

$ws2.Cells.Item(4,1).Value = $png2
$ws2.Cells.Item(4,2).Value = ".png"
$ws2.Cells.Item(4,3).Value = $status
$ws2.Cells.Item(4,4).Value = $pngTargetZh2
$ws2.Cells.Item(4,5).Value = $zhHandoffDate
$ws2.Cells.Item(4,5).NumberFormat = $dateFmt
$ws2.Cells.Item(4,8).Value = $epoch
$ws2.Cells.Item(4,8).NumberFormat = $dateFmt
$ws2.Cells.Item(4,10).Value = "IsDependency"
$ws2.Cells.Item(4,11).Value = "e2e\$newMd"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "$mdBase/$newMd", "", "", $newMd) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "$zhBase/$newXlfZh", "", "", $newXlfZh) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "$mdBase/$png1", "", "", $png1) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D3"), "$zhBase/$pngTargetZh1", "", "", $pngTargetZh1) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "$mdBase/$png2", "", "", $png2) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D4"), "$zhBase/$pngTargetZh2", "", "", $pngTargetZh2) | Out-Null

# =====================================================================
# Sheet 3: "de-de"
# =====================================================================
$ws3 = $wb.Worksheets.Item(3)
$ws3.Hyperlinks.Delete()

$ws3.Cells.Item(2,1).Value = $newMd
$ws3.Cells.Item(2,4).Value = $newXlfDe
$ws3.Cells.Item(2,5).Value = $deHandoffDate
$ws3.Cells.Item(2,5).NumberFormat = $dateFmt

$ws3.Cells.Item(3,1).Value = $png1
$ws3.Cells.Item(3,2).Value = ".png"
$ws3.Cells.Item(3,3).Value = $status
$ws3.Cells.Item(3,4).Value = $pngTargetDe1
$ws3.Cells.Item(3,5).Value = $deHandoffDate
$ws3.Cells.Item(3,5).NumberFormat = $dateFmt
$ws3.Cells.Item(3,8).Value = $epoch
$ws3.Cells.Item(3,8).NumberFormat = $dateFmt
$ws3.Cells.Item(3,10).Value = "IsDependency"
$ws3.Cells.Item(3,11).Value = "e2e\$newMd"

$ws3.Cells.Item(4,1).Value = $png2
$ws3.Cells.Item(4,2).Value = ".png"
$ws3.Cells.Item(4,3).Value = $status
$ws3.Cells.Item(4,4).Value = $pngTargetDe2
$ws3.Cells.Item(4,5).Value = $deHandoffDate
$ws3.Cells.Item(4,5).NumberFormat = $dateFmt
$ws3.Cells.Item(4,8).Value = $epoch
$ws3.Cells.Item(4,8).NumberFormat = $dateFmt
$ws3.Cells.Item(4,10).Value = "IsDependency"
$ws3.Cells.Item(4,11).Value = "e2e\$newMd"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "$mdBase/$newMd", "", "", $newMd) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "$deBase/$newXlfDe", "", "", $newXlfDe) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "$mdBase/$png1", "", "", $png1) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D3"), "$deBase/$pngTargetDe1", "", "", $pngTargetDe1) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "$mdBase/$png2", "", "", $png2) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D4"), "$deBase/$pngTargetDe2", "", "", $pngTargetDe2) | Out-Null

Write-Host "Overview hyperlinks: $($ws1.Hyperlinks.Count)"
Write-Host "zh-cn hyperlinks: $($ws2.Hyperlinks.Count)"
Write-Host "de-de hyperlinks: $($ws3.Hyperlinks.Count)"
